$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New rows of data (BLS + Census Bureau sources).
# Column-by-column so the shared-string table is populated in the same
# order the source workbook has it in (B filled before A for the BLS
# block, A filled before B for the Census block).
# ---------------------------------------------------------------------------

$blsUrls = @(
    "https://www.bls.gov/lau/laucnty13.xlsx",
    "https://www.bls.gov/lau/laucnty14.xlsx",
    "https://www.bls.gov/lau/laucnty15.xlsx",
    "https://www.bls.gov/lau/laucnty16.xlsx",
    "https://www.bls.gov/lau/laucnty17.xlsx",
    "https://www.bls.gov/lau/laucnty18.xlsx"
)
$blsLabels = @(
    "bureau of labor statistics_employment/unemployment by state_2013",
    "bureau of labor statistics_employment/unemployment by state_2014",
    "bureau of labor statistics_employment/unemployment by state_2015",
    "bureau of labor statistics_employment/unemployment by state_2016",
    "bureau of labor statistics_employment/unemployment by state_2017",
    "bureau of labor statistics_employment/unemployment by state_2018"
)

for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item(4 + $i, 2).Value = $blsUrls[$i]
}
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item(4 + $i, 1).Value = $blsLabels[$i]
}
for ($i = 0; $i -lt 6; $i++) {
    $ws.Rows.Item(4 + $i).RowHeight = 16
}

$censusLabels = @(
    "Census Bureau_2013",
    "Census Bureau_2014",
    "Census Bureau_2015",
    "Census Bureau_2016",
    "Census Bureau_2017",
    "Census Bureau_2018"
)
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item(10 + $i, 1).Value = $censusLabels[$i]
}

# Rows 10-13: Census Bureau links that carry live hyperlinks (Hyperlink style, vertical center)
$ws.Hyperlinks.Add($ws.Cells.Item(10, 2), "https://www2.census.gov/geo/tiger/GENZ2013/shp/cb_2013_us_county_20m.zip", $null, $null, "https://www2.census.gov/geo/tiger/GENZ2013/cb_2013_us_county_20m.zip")
$ws.Cells.Item(10, 2).Value = "https://www2.census.gov/geo/tiger/GENZ2013/shp/cb_2013_us_county_20m.zip"

$ws.Hyperlinks.Add($ws.Cells.Item(11, 2), "https://www2.census.gov/geo/tiger/GENZ2014/shp/cb_2014_us_county_20m.zip")
$ws.Hyperlinks.Add($ws.Cells.Item(12, 2), "https://www2.census.gov/geo/tiger/GENZ2015/shp/cb_2015_us_county_20m.zip")
$ws.Hyperlinks.Add($ws.Cells.Item(13, 2), "https://www2.census.gov/geo/tiger/GENZ2016/shp/cb_2016_us_county_20m.zip")

# Rows 14-15: Census Bureau links WITHOUT a live hyperlink (plain styled text)
$ws.Cells.Item(14, 2).Value = "https://www2.census.gov/geo/tiger/GENZ2017/shp/cb_2017_us_county_20m.zip"
$ws.Cells.Item(15, 2).Value = "https://www2.census.gov/geo/tiger/GENZ2018/shp/cb_2018_us_county_20m.zip"

# ---------------------------------------------------------------------------
# Fonts / styling, built on scratch cells far off-grid then copied in via
# PasteSpecial so the engine doesn't leave a trail of intermediate styles on
# every single cell it touches.
# ---------------------------------------------------------------------------

# Style used by B4:B9 -> Times New Roman 12, theme color 1, vertical centered
$tmp1 = $ws.Range("Z1")
$tmp1.Value = "x"
$tmp1.Font.Name = "Times New Roman"
$tmp1.Font.Size = 12
$tmp1.Font.ThemeColor = 1
$tmp1.Font.Family = 1
$tmp1.VerticalAlignment = -4108
$tmp1.Copy()
$ws.Range("B4:B9").PasteSpecial(-4122)

# Style used by B10:B13 -> Hyperlink style, vertical centered
$tmp2 = $ws.Range("Z2")
$tmp2.Value = "x"
$tmp2.Style = "Hyperlink"
$tmp2.VerticalAlignment = -4108
$tmp2.Copy()
$ws.Range("B10:B13").PasteSpecial(-4122)

# Style used by B14 -> Times New Roman 11, theme color 1, vertical centered
$tmp3 = $ws.Range("Z3")
$tmp3.Value = "x"
$tmp3.Font.Name = "Times New Roman"
$tmp3.Font.Size = 11
$tmp3.Font.ThemeColor = 1
$tmp3.Font.Family = 1
$tmp3.VerticalAlignment = -4108
$tmp3.Copy()
$ws.Range("B14").PasteSpecial(-4122)

# Style used by B15 -> Times New Roman 11, automatic color, vertical centered
$tmp4 = $ws.Range("Z4")
$tmp4.Value = "x"
$tmp4.Font.Name = "Times New Roman"
$tmp4.Font.Size = 11
$tmp4.Font.Family = 1
$tmp4.VerticalAlignment = -4108
$tmp4.Copy()
$ws.Range("B15").PasteSpecial(-4122)

$ws.Range("Z1:Z4").Clear()

# ---------------------------------------------------------------------------
# Column widths: column A was resized manually (custom width, no longer
# "best fit"); column B keeps its best-fit flag untouched.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 54

# ---------------------------------------------------------------------------
# Final selection, matching where the author's cursor was left.
# ---------------------------------------------------------------------------
$ws.Range("B19").Select()

$wb.Application.CutCopyMode = $false
